$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-07 12:43:17"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
